# Daily attendance processing - 2025-10-25 09:19:53
# Reorders the "Recorded By" list (column G) on the active sheet so that the
# last contributor in the comma-separated list is moved to the front, unless
# the last contributor is exactly "System" (in which case the cell is left
# untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G is the 7th column ("Recorded By")
$col = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value()

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ",\s*"
    if ($parts.Count -lt 2) { continue }

    $lastPart = $parts[$parts.Count - 1]
    if ($lastPart -eq "System") { continue }

    $rotated = @($lastPart) + $parts[0..($parts.Count - 2)]
    $newVal = [string]::Join(", ", $rotated)

    $cell.Value = $newVal
}
